$d = $word.ActiveDocument

# Locate the unique "a => Category" text (inside the "Constraint" paragraph)
# and collapse the found range to its end point, i.e. right after "Category"
# and before the tab that used to directly follow it.
$rng = $d.Content
$found = $rng.Find.Execute("a => Category", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'a => Category' text"
}
$rng.Collapse(0)

# Insert the new literal text "<a>" right after "Category" (it inherits the
# bold formatting of "Category" for now - the trailing tab that used to sit
# right after "Category" is pushed after the newly-inserted text, still bold).
$rng.InsertAfter("<a>")

# Re-find the freshly-inserted "<a>" text and strip its bold formatting so it
# becomes its own, non-bold run - this naturally splits the run into
# "Category" (bold) / "<a>" (non-bold) / <tab> (still bold, inherited from
# the original "Category" run) / <the remaining non-bold tabs> (untouched).
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("<a>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find inserted '<a>' text"
}
$rng2.Font.Bold = 0
